$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-12-19 Friday" "2025-12-20 Saturday"

Replace-Text "662×2=" "793×6="
Replace-Text "251×6=" "745×6="
Replace-Text "409×6=" "265×8="
Replace-Text "477×2=" "745×6="
Replace-Text "156×7=" "261×9="

Replace-Text "787×9=" "790×4="
Replace-Text "301×4=" "804×8="
Replace-Text "140×4=" "519×4="
Replace-Text "468×5=" "380×3="
Replace-Text "286×4=" "447×5="

Replace-Text "151×7=" "321×7="
Replace-Text "186×4=" "772×7="
Replace-Text "790×3=" "163×4="
Replace-Text "596×3=" "436×8="
Replace-Text "737×3=" "914×7="

Replace-Text "385×9=" "704×4="
Replace-Text "589×3=" "351×4="
Replace-Text "471×7=" "840×3="
Replace-Text "679×7=" "840×3="
Replace-Text "478×7=" "408×4="

Replace-Text "957×8=" "920×7="
Replace-Text "169×4=" "392×5="
Replace-Text "758×3=" "194×3="
Replace-Text "699×9=" "323×5="
Replace-Text "598×3=" "470×8="
